$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading-percent values for B2:O25 (case with 380 kV)
$data = @(
    @(12.10079634308016, 7.806623520735316, 5.532385992127469, 11.58203427379024, 26.63433572022007, 0, 7.344005520526261, 23.91549173517068, 0, 9.43443779019492, 9.766658145443836, 0, 19.15499729034599, 23.90746244434943),
    @(11.84286945263464, 7.771582587444176, 5.479148507636947, 11.59103829459997, 26.66347591426682, 0, 7.344005520526261, 23.99452504884335, 0, 9.252926683634863, 9.751873418653956, 0, 19.21249427670219, 23.97201195760582),
    @(11.68376696587784, 7.749725850356137, 5.445647036927164, 11.5986168572705, 26.68775224722884, 0, 7.344005520526261, 24.04722864910903, 0, 9.141145780801949, 9.744480418479366, 0, 19.24945974440285, 24.01630033891847),
    @(11.618840151422, 7.74073552356916, 5.431795931979108, 11.60222110865904, 26.69924887635836, 0, 7.344005520526261, 24.06975541000674, 0, 9.095571915532995, 9.741893993362991, 0, 19.2649426252255, 24.0355170875302),
    @(11.60805637304841, 7.739237709746671, 5.429484112745595, 11.60285076176948, 26.70125469272139, 0, 7.344005520526261, 24.07355933556625, 0, 9.088004884535557, 9.74149033632621, 0, 19.26753889712439, 24.03877856229977),
    @(11.6828915842827, 7.749604938791737, 5.445461034353205, 11.59866337605382, 26.68790080346021, 0, 7.344005520526261, 24.04752820394362, 0, 9.1405311637666, 9.744443807813489, 0, 19.24966685321435, 24.01655477218276),
    @(12.01207071404883, 7.794613675786078, 5.514201840355757, 11.58471366770203, 26.6430575122456, 0, 7.344005520526261, 23.94187511864748, 0, 9.371958097378004, 9.76121232099957, 0, 19.17447801559574, 23.92875210476652),
    @(12.64781604723911, 7.880065451616522, 5.642310592605134, 11.573600805374, 26.60581671573269, 0, 7.344005520526261, 23.7678609001782, 0, 9.820533568007528, 9.807338077962722, 0, 19.04016478107063, 23.79357204852326),
    @(13.10391677660665, 7.941001037241257, 5.732044923334734, 11.57529813437912, 26.60937968931758, 0, 7.344005520526261, 23.660274766169, 0, 10.14353965717135, 9.849114170829566, 0, 18.94941055842197, 23.7168987324118),
    @(13.30804777275081, 7.968293261224357, 5.771852039404568, 11.57819981610823, 26.61770420217657, 0, 7.344005520526261, 23.61573911751124, 0, 10.28839555327024, 9.869786500056872, 0, 18.90982808009816, 23.686950889831),
    @(13.38479004859968, 7.978564171255298, 5.786775071971398, 11.57960350616724, 26.62181816394565, 0, 7.344005520526261, 23.59950889257083, 0, 10.34289794469986, 9.877850148114723, 0, 18.89508272776646, 23.67632060275335),
    @(13.36828818258235, 7.976355036953383, 5.783567931326285, 11.57928765979964, 26.62088942460478, 0, 7.344005520526261, 23.60297613003968, 0, 10.33117630523471, 9.876103096783426, 0, 18.89824758291795, 23.67857841651375),
    @(13.31437299903581, 7.969139553743609, 5.773082823635057, 11.57830919649056, 26.61802340357041, 0, 7.344005520526261, 23.61439112654339, 0, 10.2928868279459, 9.870445207883723, 0, 18.90861009367865, 23.68606208768932),
    @(13.28127360128539, 7.964711443115448, 5.766640564892798, 11.57774952172591, 26.6163930353564, 0, 7.344005520526261, 23.62146579319922, 0, 10.26938616094696, 9.867010119700675, 0, 18.91498913221542, 23.69073858887814),
    @(13.09050230327988, 7.939208600958931, 5.729422582481519, 11.57515124425612, 26.60897037680436, 0, 7.344005520526261, 23.66327401007952, 0, 10.13402658029062, 9.847796383836684, 0, 18.95203152371473, 23.71895523006406),
    @(12.97255638267753, 7.923452068694059, 5.706327144515806, 11.57410181903316, 26.60613238297502, 0, 7.344005520526261, 23.69005108738343, 0, 10.05041677820501, 9.836433597326931, 0, 18.97519099888724, 23.73752914068237),
    @(12.90440576127786, 7.914349225857409, 5.692948174550048, 11.57369875327263, 26.60513121720889, 0, 7.344005520526261, 23.70586719831692, 0, 10.00213366866289, 9.830055320117063, 0, 18.9886719876345, 23.7486764515868),
    @(12.88128006603409, 7.911260350286736, 5.688402118927251, 11.57359676227045, 26.60490070254715, 0, 7.344005520526261, 23.71129344923836, 0, 9.985754333432983, 9.82792289047288, 0, 18.99326397869739, 23.75253040790435),
    @(12.98514470506759, 7.925133547471465, 5.708795573271603, 11.57419278657865, 26.60636918260657, 0, 7.344005520526261, 23.68715770100502, 0, 10.05933754245596, 9.837626936683936, 0, 18.97270905362506, 23.73550387549554),
    @(13.33022490981783, 7.971260675105239, 5.776166695240981, 11.57858833195216, 26.61883914704019, 0, 7.344005520526261, 23.61102103924868, 0, 10.30414331252327, 9.8721007123519, 0, 18.90555976843399, 23.68384466560804),
    @(13.55247251133379, 8.001032829154529, 5.819314777366543, 11.58323743435724, 26.63259257871167, 0, 7.344005520526261, 23.56495984128352, 0, 10.46206784295274, 9.896001682267167, 0, 18.86309364759257, 23.65422306656907),
    @(13.43417882414663, 7.985177989911685, 5.796368331222733, 11.58059407401657, 26.62474035947557, 0, 7.344005520526261, 23.58920486294764, 0, 10.37798639862415, 9.883121417262997, 0, 18.88562905697027, 23.66965344220571),
    @(12.97945458376315, 7.92437348841417, 5.707679911113568, 11.57415103626326, 26.60626016156971, 0, 7.344005520526261, 23.68846448775505, 0, 10.05530512974312, 9.837086946882662, 0, 18.97383062297893, 23.73641803731202),
    @(12.47741139755568, 7.857263008752606, 5.608402410629573, 11.57487169976249, 26.61045916878746, 0, 7.344005520526261, 23.81138094591088, 0, 9.70010032817901, 9.793460124805856, 0, 19.07510250535228, 23.82617115558732)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item(2 + $i, 2 + $j).Value = $rowValues[$j]
    }
}

